# Add two new Mac-Addresses (10 new detail rows) to the reg_center_machine_device
# master data sheet, mirroring the existing row layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (row, regcntr_id, machine_id, device_id)
$newRows = @(
    @(147, 10001, 10030, 3000166),
    @(148, 10001, 10030, 3000167),
    @(149, 10001, 10030, 3000168),
    @(150, 10001, 10030, 3000169),
    @(151, 10001, 10030, 3000170),
    @(152, 10001, 10031, 3000171),
    @(153, 10001, 10031, 3000172),
    @(154, 10001, 10031, 3000173),
    @(155, 10001, 10031, 3000174),
    @(156, 10001, 10031, 3000175)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = "eng"
    $ws.Cells.Item($rowNum, 5).Value = $true
    $ws.Cells.Item($rowNum, 6).Value = "superadmin"
    $ws.Cells.Item($rowNum, 7).Value = "now()"
}

# Scroll the view down to reflect the newly added rows and move the active
# selection to the first empty cell after the data.
$ws.Application.ActiveWindow.ScrollRow = 144
[void]$ws.Range("H149").Select()
